{"js": "// The contribution-statement paragraph is rewritten in its entirety (new\n// wording about NISAR/UAVSAR and global SWE estimation). There is a single\n// paragraph in the document body, so grab it and replace its full text in\n// one shot; \"Replace\" keeps the paragraph (and its justified alignment)\n// while swapping out every run inside it for one run with the new text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\nconst newText =\n  \"Quantifying snow water equivalent (SWE) on a global scale has posed a long-standing challenge. \" +\n  \"SWE can be calculated if the snowpack depth and density measurements are known. \" +\n  \"However, these products are only available in selected locations worldwide. \" +\n  \"The forthcoming NASA-ISRO Synthetic Aperture Radar (NISAR) mission offers a new horizon for snow monitoring, with its global coverage and frequent revisit cycles. \" +\n  \"The data from this mission can be used to estimate changes in snow depth and SWE. \" +\n  \"Our work uses NISAR-like data from the NASA JPL Uninhabited Aerial Vehicle Synthetic Aperture Radar (UAVSAR) sensor to estimate total snow depth using machine learning (ML) algorithms. \" +\n  \"We demonstrate the potential of combining advanced radar technology with machine learning algorithms to produce snow depth maps. \" +\n  \"Our findings are a step towards developing a global snow depth prediction system that will provide valuable information for water resource management, flood forecasting, and avalanche hazard assessment, provided that accurate and representative training data is available. \" +\n  \"By showcasing the effectiveness of UAVSAR data in snow depth estimation, our research highlights a path forward for snow monitoring in anticipation of the capabilities that the NISAR mission will further expand upon.\";\n\nparagraph.insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# The contribution-statement paragraph is rewritten in its entirety (new\n# wording about NISAR/UAVSAR and global SWE estimation). There is a single\n# paragraph in the document body. Paragraph.Range.Text includes the\n# trailing paragraph mark, so build a Range over just the paragraph's\n# characters (excluding that mark) and overwrite its Text in one shot -\n# this preserves the paragraph mark and paragraph formatting (justified\n# alignment) while swapping out every run for one run with the new text.\n$d = $word.ActiveDocument\n$p = $d.Paragraphs.First\n$charCount = $p.Range.Text.Length - 1\n\n$newText = \"Quantifying snow water equivalent (SWE) on a global scale has posed a long-standing challenge. \" + `\n  \"SWE can be calculated if the snowpack depth and density measurements are known. \" + `\n  \"However, these products are only available in selected locations worldwide. \" + `\n  \"The forthcoming NASA-ISRO Synthetic Aperture Radar (NISAR) mission offers a new horizon for snow monitoring, with its global coverage and frequent revisit cycles. \" + `\n  \"The data from this mission can be used to estimate changes in snow depth and SWE. \" + `\n  \"Our work uses NISAR-like data from the NASA JPL Uninhabited Aerial Vehicle Synthetic Aperture Radar (UAVSAR) sensor to estimate total snow depth using machine learning (ML) algorithms. \" + `\n  \"We demonstrate the potential of combining advanced radar technology with machine learning algorithms to produce snow depth maps. \" + `\n  \"Our findings are a step towards developing a global snow depth prediction system that will provide valuable information for water resource management, flood forecasting, and avalanche hazard assessment, provided that accurate and representative training data is available. \" + `\n  \"By showcasing the effectiveness of UAVSAR data in snow depth estimation, our research highlights a path forward for snow monitoring in anticipation of the capabilities that the NISAR mission will further expand upon.\"\n\n$target = $d.Range(0, $charCount)\n$target.Text = $newText\n"}
